$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 20250318
$ws.Range("C1").Value = "Ruta de ejemplo PTO MONTT"

# Row 3
$ws.Range("A3").Value = 20250318
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 9999
$ws.Range("E3").Value = "ACHS VALDIVIA"
$ws.Range("F3").Value = "Beauchef Nº705, Valdivia"
$ws.Range("G3").Value = "Valdivia"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = $null

# Row 4
$ws.Range("A4").Value = 20250318
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 60
$ws.Range("D4").Value = "77.880.747-5"
$ws.Range("E4").Value = "Salud Humana Pablo Contreras"
$ws.Range("F4").Value = "beauchef 925  ( nueva direccion)"
$ws.Range("G4").Value = "Valdivia"
$ws.Range("H4").Value = 966796589
$ws.Range("I4").Value = "-"
$ws.Range("J4").Value = $null

# Row 5
$ws.Range("A5").Value = 20250318
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 60
$ws.Range("D5").Value = "77.624.403-1"
$ws.Range("E5").Value = "serv para salud humana Gloria Jofré EIRL"
$ws.Range("F5").Value = "García Reyes 686"
$ws.Range("G5").Value = "Valdivia"
$ws.Range("H5").Value = 952541245
$ws.Range("I5").Value = "DEUDA"
$ws.Range("J5").Value = $null

# Row 6
$ws.Range("A6").Value = 20250318
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "120"
$ws.Range("D6").Value = 9999
$ws.Range("E6").Value = "DOCTOR SIMI"
$ws.Range("F6").Value = "mall valdivia"
$ws.Range("G6").Value = "valdivia"
$ws.Range("H6").Value = "-"
$ws.Range("I6").Value = "SERVICIO ESPECIAL"
$ws.Range("J6").Value = $null

# Row 7
$ws.Range("A7").Value = 20250318
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 9999
$ws.Range("E7").Value = "ACHS Osorno"
$ws.Range("F7").Value = "Av Zenteno 1529"
$ws.Range("G7").Value = "Osorno"
$ws.Range("H7").Value = "-"
$ws.Range("I7").Value = $null
$ws.Range("J7").Value = $null

# Row 8
$ws.Range("A8").Value = 20250318
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = "65.062.843-8"
$ws.Range("E8").Value = "Fundacion salud y Familia"
$ws.Range("F8").Value = "Benavente 379 of 206"
$ws.Range("G8").Value = "Puerto Montt"
$ws.Range("H8").Value = 987516103
$ws.Range("I8").Value = "a las 16:00 / nesecita 100 bolsas amarillas"
$ws.Range("J8").Value = $null

# Row 9
$ws.Range("A9").Value = 20250318
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 15
$ws.Range("D9").Value = 9999
$ws.Range("E9").Value = "ACHS PUERTO MONTT"
$ws.Range("F9").Value = "Ejército Nº360"
$ws.Range("G9").Value = "Puerto Montt"
$ws.Range("H9").Value = "-"
$ws.Range("I9").Value = $null
$ws.Range("J9").Value = $null

# Row 10
$ws.Range("A10").Value = 20250318
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 9999
$ws.Range("E10").Value = "EVEREST PTO MONTT"
$ws.Range("F10").Value = "Juan Soler Manfredini L 131"
$ws.Range("G10").Value = "Puerto Montt"
$ws.Range("H10").Value = "-"
$ws.Range("I10").Value = "8:30 a 13 y 14 a 19"
$ws.Range("J10").Value = $null

# Row 11
$ws.Range("A11").Value = 20250318
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = "65.062.843-8"
$ws.Range("E11").Value = "CESFAM San Pablo"
$ws.Range("F11").Value = "Barrancas S/N mirasol"
$ws.Range("G11").Value = "Puerto Montt"
$ws.Range("H11").Value = 990579442
$ws.Range("I11").Value = "9:30 a 14 y 17:30 a 20"
$ws.Range("J11").Value = $null

# Row 12
$ws.Range("A12").Value = 20250318
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 60
$ws.Range("D12").Value = "77.628.398-3"
$ws.Range("E12").Value = "CLINICA ODONTOLOGICA ACDENT SPA"
$ws.Range("F12").Value = "Antonio Varas 55 of 509"
$ws.Range("G12").Value = "Puerto Montt"
$ws.Range("H12").Value = "964988313/967120480"
$ws.Range("I12").Value = "10 a 13 y 15 a 17:30"
$ws.Range("J12").Value = $null

# Row 13
$ws.Range("A13").Value = 20250318
$ws.Range("B13").Value = 11
$ws.Range("C13").Value = 60
$ws.Range("D13").Value = "77.491.661-k"
$ws.Range("E13").Value = "Emuna soluciones medicas spa"
$ws.Range("F13").Value = "Benavente 840,Of 605"
$ws.Range("G13").Value = "Puerto Montt"
$ws.Range("H13").Value = "962826997/951701572"
$ws.Range("I13").Value = "9 a 13 y 14 a 16:00"
$ws.Range("J13").Value = $null

# Row 14
$ws.Range("A14").Value = 20250318
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 60
$ws.Range("D14").Value = "76.058.615-3"
$ws.Range("E14").Value = "Fernandez y Haussmann odontologia"
$ws.Range("F14").Value = "Quillota 175 of 1110"
$ws.Range("G14").Value = "Puerto Montt"
$ws.Range("H14").Value = "652263740/934533835"
$ws.Range("I14").Value = "10 a 18"
$ws.Range("J14").Value = $null

# Row 15
$ws.Range("A15").Value = 20250318
$ws.Range("B15").Value = 13
$ws.Range("C15").Value = 60
$ws.Range("D15").Value = "78.034.130-0"
$ws.Range("E15").Value = "KARUS"
$ws.Range("F15").Value = "Juan soler manfredini 41 of 1803"
$ws.Range("G15").Value = "Puerto Montt"
$ws.Range("H15").Value = "952296251/957121954"
$ws.Range("I15").Value = "cliente nuevo 1C Y 1B /11 a  14"
$ws.Range("J15").Value = $null

# Row 16
$ws.Range("A16").Value = 20250318
$ws.Range("B16").Value = 14
$ws.Range("C16").Value = 60
$ws.Range("D16").Value = "77.962.205-3"
$ws.Range("E16").Value = "we face and body spa."
$ws.Range("F16").Value = "O´higgins 167 of 609"
$ws.Range("G16").Value = "Puerto Montt"
$ws.Range("H16").Value = 972078066
$ws.Range("I16").Value = "10 a 19"
$ws.Range("J16").Value = $null

# Row 17
$ws.Range("A17").Value = 20250318
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = 60
$ws.Range("D17").Value = "77.141.688-8"
$ws.Range("E17").Value = "Servicios Sociales Sur Spa"
$ws.Range("F17").Value = "Francisco Bilbao 255"
$ws.Range("G17").Value = "Puerto Montt"
$ws.Range("H17").Value = 990847901
$ws.Range("I17").Value = "9 a 13 y de 14 a 17"
$ws.Range("J17").Value = $null

# Row 18
$ws.Range("A18").Value = 20250318
$ws.Range("B18").Value = 16
$ws.Range("C18").Value = 60
$ws.Range("D18").Value = "78.074.427-8"
$ws.Range("E18").Value = "serv. Odontologicos Loreto sanz"
$ws.Range("F18").Value = "Santa rosa 575 of 27"
$ws.Range("G18").Value = "Puerto Varas"
$ws.Range("H18").Value = 982188029
$ws.Range("I18").Value = "10 a 18"
$ws.Range("J18").Value = $null

